$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -887
$ws.Range("H5").Value = 211.85715
$ws.Range("I5").Value = 233.5
$ws.Range("K5").Value = 233.5
$ws.Range("M5").Value = -118.5
$ws.Range("H17").Value = 345334.56
$ws.Range("J17").Value = 345334.56
$ws.Range("L17").Value = 1036003.68
$ws.Range("N17").Value = -1036339.68
$ws.Range("H43").Value = 1101.6666
$ws.Range("I43").Value = 1077.5
$ws.Range("K43").Value = 1077.5
$ws.Range("M43").Value = -1008.5
$ws.Range("H53").Value = 101316.4
$ws.Range("I53").Value = 322
$ws.Range("J53").Value = 202310.8
$ws.Range("K53").Value = 322
$ws.Range("L53").Value = 202310.8
$ws.Range("M53").Value = 315
$ws.Range("N53").Value = -203584.8
$ws.Range("H113").Value = 83337660
$ws.Range("I113").Value = 250003000
$ws.Range("K113").Value = 250003000
$ws.Range("M113").Value = -249999746
$ws.Range("H116").Value = 58451668
$ws.Range("I116").Value = 35874290
$ws.Range("K116").Value = 35874290
$ws.Range("M116").Value = -35870848
$ws.Range("H131").Value = 7405.067
$ws.Range("I131").Value = 852.0769
$ws.Range("K131").Value = 2556.2307
$ws.Range("M131").Value = 2483.7693
$ws.Range("H132").Value = 3553.1892
$ws.Range("I132").Value = 3199.4517
$ws.Range("K132").Value = 9598.355100000001
$ws.Range("M132").Value = -7068.355100000001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1440.1428
$ws.Range("I2").Value = 1524.091
$ws.Range("K2").Value = 1524.091
$ws.Range("M2").Value = -1411.091
$ws.Range("H32").Value = 3902.6155
$ws.Range("I32").Value = 2426.1936
$ws.Range("K32").Value = 2426.1936
$ws.Range("M32").Value = -2139.1936
$ws.Range("H61").Value = 37040410
$ws.Range("I61").Value = 47621450
$ws.Range("K61").Value = 47621450
$ws.Range("M61").Value = -47621238
$ws.Range("H116").Value = 1440.1428
$ws.Range("I116").Value = 1524.091
$ws.Range("K116").Value = 1524.091
$ws.Range("M116").Value = 769.9090000000001
$ws.Range("H132").Value = 30305046
$ws.Range("I132").Value = 37038904
$ws.Range("K132").Value = 111116712
$ws.Range("M132").Value = -111114182
$ws.Range("H135").Value = 46880.5
$ws.Range("J135").Value = 46880.5
$ws.Range("L135").Value = 46880.5
$ws.Range("N135").Value = -57020.5
$ws.Range("H136").Value = 37040410
$ws.Range("I136").Value = 47621450
$ws.Range("K136").Value = 142864350
$ws.Range("M136").Value = -142861800

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1440.1428
$ws.Range("I3").Value = 1524.091
$ws.Range("K3").Value = 1524.091
$ws.Range("M3").Value = -1410.091

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 1621
$ws.Range("I33").Value = 1621
$ws.Range("K33").Value = 1621
$ws.Range("M33").Value = -1242
$ws.Range("H99").Value = 2789.6428
$ws.Range("I99").Value = 2416.6667
$ws.Range("K99").Value = 2416.6667
$ws.Range("M99").Value = -918.6667000000002
$ws.Range("H122").Value = 2284.8823
$ws.Range("I122").Value = 2254.25
$ws.Range("J122").Value = 2358.4
$ws.Range("K122").Value = 6762.75
$ws.Range("L122").Value = 7075.200000000001
$ws.Range("M122").Value = -4312.75
$ws.Range("N122").Value = -11975.2
$ws.Range("H126").Value = 2789.6428
$ws.Range("I126").Value = 2416.6667
$ws.Range("K126").Value = 7250.000100000001
$ws.Range("M126").Value = -4780.000100000001
$ws.Range("H132").Value = 3550.0527
$ws.Range("I132").Value = 3469.5557
$ws.Range("K132").Value = 10408.6671
$ws.Range("M132").Value = -7878.667099999999
$ws.Range("H134").Value = 2263.3333
$ws.Range("I134").Value = 1717.1428
$ws.Range("K134").Value = 5151.428400000001
$ws.Range("M134").Value = -2616.428400000001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8338863
$ws.Range("I4").Value = 12060950
$ws.Range("K4").Value = 36182850
$ws.Range("M4").Value = -36182738
$ws.Range("H68").Value = 547.125
$ws.Range("J68").Value = 800
$ws.Range("L68").Value = 2400
$ws.Range("N68").Value = -4022
$ws.Range("H71").Value = 547.125
$ws.Range("J71").Value = 800
$ws.Range("L71").Value = 7200
$ws.Range("N71").Value = -15312
$ws.Range("H74").Value = 5338.3335
$ws.Range("J74").Value = 7507.5
$ws.Range("L74").Value = 22522.5
$ws.Range("N74").Value = -24644.5
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 5338.3335
$ws.Range("J77").Value = 7507.5
$ws.Range("L77").Value = 67567.5
$ws.Range("N77").Value = -78175.5
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("H81").Value = 7449.5864
$ws.Range("J81").Value = 8494.559999999999
$ws.Range("L81").Value = 25483.68
$ws.Range("N81").Value = -27729.68
$ws.Range("H84").Value = 7449.5864
$ws.Range("J84").Value = 8494.559999999999
$ws.Range("L84").Value = 76451.03999999999
$ws.Range("N84").Value = -87683.03999999999
$ws.Range("H87").Value = 1833
$ws.Range("I87").Value = 1833
$ws.Range("K87").Value = 5499
$ws.Range("M87").Value = -4251
$ws.Range("H90").Value = 1833
$ws.Range("I90").Value = 1833
$ws.Range("K90").Value = 16497
$ws.Range("M90").Value = -10257
$ws.Range("H107").Value = 909.2778
$ws.Range("I107").Value = 889.26666
$ws.Range("J107").Value = 1009.3333
$ws.Range("K107").Value = 2667.79998
$ws.Range("L107").Value = 3027.9999
$ws.Range("M107").Value = -747.7999799999998
$ws.Range("N107").Value = -6867.9999
$ws.Range("H131").Value = 12537.588
$ws.Range("J131").Value = 13972.6
$ws.Range("L131").Value = 41917.8
$ws.Range("N131").Value = -51997.8

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 13
$ws.Range("N2").ClearContents()
$ws.Range("H43").Value = 2768.8667
$ws.Range("I43").Value = 2768.8667
$ws.Range("K43").Value = 2768.8667
$ws.Range("M43").Value = -2617.8667
$ws.Range("H102").Value = 2090.4736
$ws.Range("I102").Value = 1607.44
$ws.Range("K102").Value = 1607.44
$ws.Range("M102").Value = 14.55999999999995
$ws.Range("H132").Value = 4563.7144
$ws.Range("I132").Value = 4267.8823
$ws.Range("J132").Value = 5020.909
$ws.Range("K132").Value = 12803.6469
$ws.Range("L132").Value = 15062.727
$ws.Range("M132").Value = -10273.6469
$ws.Range("N132").Value = -20122.727

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 385199.8
$ws.Range("J4").Value = 231499.75
$ws.Range("L4").Value = 231499.75
$ws.Range("N4").Value = -231725.75
$ws.Range("H28").Value = 385199.8
$ws.Range("J28").Value = 231499.75
$ws.Range("L28").Value = 231499.75
$ws.Range("N28").Value = -231963.75
$ws.Range("H37").Value = 385199.8
$ws.Range("J37").Value = 231499.75
$ws.Range("L37").Value = 231499.75
$ws.Range("N37").Value = -231713.75
$ws.Range("H55").Value = 420.54544
$ws.Range("I55").Value = 405.8
$ws.Range("J55").Value = 452.14285
$ws.Range("K55").Value = 405.8
$ws.Range("L55").Value = 452.14285
$ws.Range("M55").Value = -232.8
$ws.Range("N55").Value = -798.14285
$ws.Range("H122").Value = 3632.3076
$ws.Range("I122").Value = 2761.1765
$ws.Range("K122").Value = 8283.529500000001
$ws.Range("M122").Value = -5833.529500000001
$ws.Range("H132").Value = 5120.421
$ws.Range("I132").Value = 3350.0908
$ws.Range("K132").Value = 10050.2724
$ws.Range("M132").Value = -7520.2724

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H63").Value = 29995
$ws.Range("J63").Value = 29995
$ws.Range("L63").Value = 29995
$ws.Range("N63").Value = -31243
$ws.Range("H66").Value = 29995
$ws.Range("J66").Value = 29995
$ws.Range("L66").Value = 89985
$ws.Range("N66").Value = -96225
$ws.Range("H122").Value = 2507.1538
$ws.Range("I122").Value = 2417.6365
$ws.Range("K122").Value = 7252.9095
$ws.Range("M122").Value = -4802.9095
$ws.Range("H125").Value = 54999
$ws.Range("J125").Value = 54999
$ws.Range("L125").Value = 54999
$ws.Range("N125").Value = -64839
